$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 900858550
$ws.Range("C2").Value = "BIBO SOLUTIONS"
$ws.Range("D2").Value = "SAS"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = "CA"
$ws.Range("G2").Value = "83645289326"
$ws.Range("H2").Value = 2359386

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 900654100
$ws.Range("C3").Value = "CIMAZ S.A.S"
$ws.Range("H3").Value = 525870

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 1143940723
$ws.Range("C4").Value = "CIMPRE"
$ws.Range("D4").Value = "SALUD OCUPACIONAL S.A.S."
$ws.Range("H4").Value = 87451

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 31322510
$ws.Range("C5").Value = "IMPATA RESTREPO DIANA CARINA"
$ws.Range("H5").Value = 84300

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 901223156
$ws.Range("C6").Value = "DIGITALTIC SAS"
$ws.Range("H6").Value = 193970

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 7215649
$ws.Range("C7").Value = "DOMINGO IGNACIO"
$ws.Range("D7").Value = "ROJAS"
$ws.Range("E7").Value = 32
$ws.Range("F7").Value = "CA"
$ws.Range("G7").Value = "24003323467"
$ws.Range("H7").Value = 97991

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 890304345
$ws.Range("C8").Value = "ELECTRICOS DEL VALLE SA"
$ws.Range("H8").Value = 172500

$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 890306372
$ws.Range("C9").Value = "ELECTRO JAPONESA S.A."
$ws.Range("H9").Value = 1819546

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 900298074
$ws.Range("C10").Value = "GVS COLOMBIA SAS"
$ws.Range("H10").Value = 9952604

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 830076882
$ws.Range("C11").Value = "Hp Financial Services Colombia LLC Sucursal Colombia"
$ws.Range("H11").Value = 4407849

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = 1143940722
$ws.Range("C12").Value = "IZC"
$ws.Range("D12").Value = "MAYORISTA SAS"
$ws.Range("H12").Value = 52092009

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = 94281756
$ws.Range("C13").Value = "JUAN CARLOS"
$ws.Range("D13").Value = "MARQUEZ SANCHEZ"
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "CA"
$ws.Range("G13").Value = "03165339508"
$ws.Range("H13").Value = 1074121

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 900892841
$ws.Range("C14").Value = "LILIUM TECNOLOGIA SAS"
$ws.Range("H14").Value = 669600

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 800035776
$ws.Range("C15").Value = "NEXSYS DE COLOMBIA SA"
$ws.Range("H15").Value = 18089916

$ws.Range("A16").Value = 1
$ws.Range("B16").Value = 830034343
$ws.Range("C16").Value = "RENTEK SAS"
$ws.Range("H16").Value = 4094318

$ws.Range("A17").Value = 1
$ws.Range("B17").Value = 900355222
$ws.Range("C17").Value = "TIENDAS TECNOPLAZA S.A.S"
$ws.Range("H17").Value = 492503

$ws.Range("A18").Value = 3
$ws.Range("B18").Value = 7685100
$ws.Range("C18").Value = "ROJAS SALAZAR WILLIAM"
$ws.Range("H18").Value = 6013778

$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 800179308
$ws.Range("C19").Value = "YAMAKI SAS"
$ws.Range("H19").Value = 4031339

$ws.Columns.Item(4).ColumnWidth = 25
